# Update loading_percent values for the 380 kV case (Case_4_106)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 17.71478051727053
$ws.Range("C2").Value = 3.918362450550199
$ws.Range("D2").Value = 13.24758217417682
$ws.Range("E2").Value = 13.13359281170308
$ws.Range("G2").Value = 60.71641640489961
$ws.Range("H2").Value = 22.34760432598437
$ws.Range("I2").Value = 35.2352596604217
$ws.Range("J2").Value = 8.061198718603105
$ws.Range("K2").Value = 15.15684806789466
$ws.Range("L2").Value = 13.027669503121
$ws.Range("M2").Value = 18.48091113528358
$ws.Range("N2").Value = 24.61537505153531

$ws.Range("B3").Value = 17.60320496475338
$ws.Range("C3").Value = 3.701879639757636
$ws.Range("D3").Value = 13.25114696125519
$ws.Range("E3").Value = 13.15536450863131
$ws.Range("G3").Value = 60.75684234159249
$ws.Range("H3").Value = 22.39061415335945
$ws.Range("I3").Value = 35.30440407070573
$ws.Range("J3").Value = 8.060591239375906
$ws.Range("K3").Value = 15.08107307220823
$ws.Range("L3").Value = 13.04398961377904
$ws.Range("M3").Value = 18.48133287078555
$ws.Range("N3").Value = 24.67460215187332

$ws.Range("B4").Value = 17.53823456229862
$ws.Range("C4").Value = 3.595258464558431
$ws.Range("D4").Value = 13.25545608998988
$ws.Range("E4").Value = 13.16992849902011
$ws.Range("G4").Value = 60.79357584414558
$ws.Range("H4").Value = 22.42006014324332
$ws.Range("I4").Value = 35.35196957328785
$ws.Range("J4").Value = 8.060236966094333
$ws.Range("K4").Value = 15.03754913401113
$ws.Range("L4").Value = 13.05561425727466
$ws.Range("M4").Value = 18.48438561187706
$ws.Range("N4").Value = 24.7128745041517

$ws.Range("B5").Value = 17.5126719495962
$ws.Range("C5").Value = 3.552756225930684
$ws.Range("D5").Value = 13.25774611291146
$ws.Range("E5").Value = 13.17616472975301
$ws.Range("G5").Value = 60.81153804162208
$ws.Range("H5").Value = 22.43282337484348
$ws.Range("I5").Value = 35.37263729333726
$ws.Range("J5").Value = 8.06009730141937
$ws.Range("K5").Value = 15.02058238216896
$ws.Range("L5").Value = 13.06075517545371
$ws.Range("M5").Value = 18.48633393125439
$ws.Range("N5").Value = 24.72895125243338

$ws.Range("B6").Value = 17.5084831294384
$ws.Range("C6").Value = 3.545664333500641
$ws.Range("D6").Value = 13.25815864401154
$ws.Range("E6").Value = 13.1772184633361
$ws.Range("G6").Value = 60.81470133111922
$ws.Range("H6").Value = 22.43498882397833
$ws.Range("I6").Value = 35.37614670435654
$ws.Range("J6").Value = 8.060074394262875
$ws.Range("K6").Value = 15.0178119562339
$ws.Range("L6").Value = 13.06163322034332
$ws.Range("M6").Value = 18.4867000336509
$ws.Range("N6").Value = 24.73164983218583

$ws.Range("B7").Value = 17.53788608743498
$ws.Range("C7").Value = 3.594687622306843
$ws.Range("D7").Value = 13.25548481083376
$ws.Range("E7").Value = 13.17001138241714
$ws.Range("G7").Value = 60.79380597405336
$ws.Range("H7").Value = 22.42022918047032
$ws.Range("I7").Value = 35.35224310596703
$ws.Range("J7").Value = 8.060235063481686
$ws.Range("K7").Value = 15.03731717977009
$ws.Range("L7").Value = 13.05568195405387
$ws.Range("M7").Value = 18.48440903373261
$ws.Range("N7").Value = 24.71308937414401

$ws.Range("B8").Value = 17.6755898068776
$ws.Range("C8").Value = 3.845271364376526
$ws.Range("D8").Value = 13.24837179168995
$ws.Range("E8").Value = 13.14085176093701
$ws.Range("G8").Value = 60.72788187972945
$ws.Range("H8").Value = 22.361803560568
$ws.Range("I8").Value = 35.25803956210999
$ws.Range("J8").Value = 8.060985372000429
$ws.Range("K8").Value = 15.13010685996911
$ws.Range("L8").Value = 13.03296403288619
$ws.Range("M8").Value = 18.48047783283536
$ws.Range("N8").Value = 24.63540143232323

$ws.Range("B9").Value = 17.97258133614159
$ws.Range("C9").Value = 4.343658567823863
$ws.Range("D9").Value = 13.25120428117574
$ws.Range("E9").Value = 13.09313712113942
$ws.Range("G9").Value = 60.69318508426353
$ws.Range("H9").Value = 22.27133954088407
$ws.Range("I9").Value = 35.11388630775974
$ws.Range("J9").Value = 8.062605530141326
$ws.Range("K9").Value = 15.33519613200596
$ws.Range("L9").Value = 13.0011216599864
$ws.Range("M9").Value = 18.49485011218497
$ws.Range("N9").Value = 24.49813777773715

$ws.Range("B10").Value = 18.20563734382365
$ws.Range("C10").Value = 4.672913801497349
$ws.Range("D10").Value = 13.26344943594063
$ws.Range("E10").Value = 13.06382146775981
$ws.Range("G10").Value = 60.72538389596781
$ws.Range("H10").Value = 22.21957877487881
$ws.Range("I10").Value = 35.03275333803433
$ws.Range("J10").Value = 8.063887108504012
$ws.Range("K10").Value = 15.49898412335085
$ws.Range("L10").Value = 12.98544602533603
$ws.Range("M10").Value = 18.51874346592655
$ws.Range("N10").Value = 24.40641890620403

$ws.Range("B11").Value = 18.31452903430927
$ws.Range("C11").Value = 4.81461683997636
$ws.Range("D11").Value = 13.27120926741477
$ws.Range("E11").Value = 13.05172479064737
$ws.Range("G11").Value = 60.75254179460727
$ws.Range("H11").Value = 22.19922452230399
$ws.Range("I11").Value = 35.00123006280054
$ws.Range("J11").Value = 8.064490070688365
$ws.Range("K11").Value = 15.57611315209197
$ws.Range("L11").Value = 12.9799840955429
$ws.Range("M11").Value = 18.53247631151154
$ws.Range("N11").Value = 24.36666197965595

$ws.Range("B12").Value = 18.35614538094997
$ws.Range("C12").Value = 4.867111283798312
$ws.Range("D12").Value = 13.27446068164802
$ws.Range("E12").Value = 13.04732174079235
$ws.Range("G12").Value = 60.76462099440128
$ws.Range("H12").Value = 22.19197583339879
$ws.Range("I12").Value = 34.99006760928479
$ws.Range("J12").Value = 8.064721274308976
$ws.Range("K12").Value = 15.60567642605488
$ws.Range("L12").Value = 12.97815513735031
$ws.Range("M12").Value = 18.53808516212131
$ws.Range("N12").Value = 24.35188884601176

$ws.Range("B13").Value = 18.34716610388527
$ws.Range("C13").Value = 4.855857544742618
$ws.Range("D13").Value = 13.27374654602173
$ws.Range("E13").Value = 13.04826212099583
$ws.Range("G13").Value = 60.76193976108407
$ws.Range("H13").Value = 22.19351655329174
$ws.Range("I13").Value = 34.99243718152992
$ws.Range("J13").Value = 8.064671352547974
$ws.Range("K13").Value = 15.59929394084399
$ws.Range("L13").Value = 12.97853840159707
$ws.Range("M13").Value = 18.53685908157672
$ws.Range("N13").Value = 24.35505797768372

$ws.Range("B14").Value = 18.31794537582796
$ws.Range("C14").Value = 4.818958963284435
$ws.Range("D14").Value = 13.27147050316492
$ws.Range("E14").Value = 13.05135899055291
$ws.Range("G14").Value = 60.75349960820722
$ws.Range("H14").Value = 22.19861896885671
$ws.Range("I14").Value = 35.00029619228818
$ws.Range("J14").Value = 8.064509034309269
$ws.Range("K14").Value = 15.57853832275669
$ws.Range("L14").Value = 12.97982883287017
$ws.Range("M14").Value = 18.53292958732908
$ws.Range("N14").Value = 24.36544093897457

$ws.Range("B15").Value = 18.3000955564985
$ws.Range("C15").Value = 4.796205654158226
$ws.Range("D15").Value = 13.2701170545707
$ws.Range("E15").Value = 13.05327903963335
$ws.Range("G15").Value = 60.74856340732423
$ws.Range("H15").Value = 22.20180412390846
$ws.Range("I15").Value = 35.0052109664741
$ws.Range("J15").Value = 8.064409984217397
$ws.Range("K15").Value = 15.56587066797045
$ws.Range("L15").Value = 12.98065041052712
$ws.Range("M15").Value = 18.5305757585351
$ws.Range("N15").Value = 24.37183749690341

$ws.Range("B16").Value = 18.19857613432179
$ws.Range("C16").Value = 4.663490269494276
$ws.Range("D16").Value = 13.26298621163724
$ws.Range("E16").Value = 13.06463690729916
$ws.Range("G16").Value = 60.72386050232106
$ws.Range("H16").Value = 22.22097320720811
$ws.Range("I16").Value = 35.03492184891089
$ws.Range("J16").Value = 8.0638481056004
$ws.Range("K16").Value = 15.49399454507742
$ws.Range("L16").Value = 12.98583650593376
$ws.Range("M16").Value = 18.51790334258678
$ws.Range("N16").Value = 24.40905661422319

$ws.Range("B17").Value = 18.13701079800773
$ws.Range("C17").Value = 4.580001634226346
$ws.Range("D17").Value = 13.25917114788586
$ws.Range("E17").Value = 13.07192162226739
$ws.Range("G17").Value = 60.71190793887023
$ws.Range("H17").Value = 22.23355039635578
$ws.Range("I17").Value = 35.05452796974387
$ws.Range("J17").Value = 8.063508527731358
$ws.Range("K17").Value = 15.45055764297277
$ws.Range("L17").Value = 12.98944503591667
$ws.Range("M17").Value = 18.51086043124619
$ws.Range("N17").Value = 24.43239240200231

$ws.Range("B18").Value = 18.10187244306796
$ws.Range("C18").Value = 4.531221070575407
$ws.Range("D18").Value = 13.25718303204406
$ws.Range("E18").Value = 13.07622826290742
$ws.Range("G18").Value = 60.70621100189308
$ws.Range("H18").Value = 22.24108490948713
$ws.Range("I18").Value = 35.06631168974904
$ws.Range("J18").Value = 8.063315090566835
$ws.Range("K18").Value = 15.42582187361105
$ws.Range("L18").Value = 12.99167770239024
$ws.Range("M18").Value = 18.50707928644745
$ws.Range("N18").Value = 24.44599966926603

$ws.Range("B19").Value = 18.09002294176848
$ws.Range("C19").Value = 4.514574409619672
$ws.Range("D19").Value = 13.25654536140824
$ws.Range("E19").Value = 13.07770646967898
$ws.Range("G19").Value = 60.70448451377057
$ws.Range("H19").Value = 22.24368756653547
$ws.Range("I19").Value = 35.07038848251426
$ws.Range("J19").Value = 8.06324991867521
$ws.Range("K19").Value = 15.41748998041321
$ws.Range("L19").Value = 12.99246065181679
$ws.Range("M19").Value = 18.50584548113255
$ws.Range("N19").Value = 24.45063867934946

$ws.Range("B20").Value = 18.14353656904526
$ws.Range("C20").Value = 4.588967844501305
$ws.Range("D20").Value = 13.25955593937153
$ws.Range("E20").Value = 13.0711340810493
$ws.Range("G20").Value = 60.71305842955044
$ws.Range("H20").Value = 22.23218043762568
$ws.Range("I20").Value = 35.05238840941604
$ws.Range("J20").Value = 8.063544481874185
$ws.Range("K20").Value = 15.45515605767412
$ws.Range("L20").Value = 12.98904464330262
$ws.Range("M20").Value = 18.51158226573994
$ws.Range("N20").Value = 24.42988911230687

$ws.Range("B21").Value = 18.32651811617269
$ws.Range("C21").Value = 4.829828623980506
$ws.Range("D21").Value = 13.27213055534034
$ws.Range("E21").Value = 13.05044454611948
$ws.Range("G21").Value = 60.75593000482429
$ws.Range("H21").Value = 22.19710780804522
$ws.Range("I21").Value = 34.99796678158386
$ws.Range("J21").Value = 8.064556633055204
$ws.Range("K21").Value = 15.58462525223628
$ws.Range("L21").Value = 12.97944331125615
$ws.Range("M21").Value = 18.53407271527113
$ws.Range("N21").Value = 24.36238356593013

$ws.Range("B22").Value = 18.44831527076828
$ws.Range("C22").Value = 4.980451340541854
$ws.Range("D22").Value = 13.28217179201452
$ws.Range("E22").Value = 13.03795827553223
$ws.Range("G22").Value = 60.79440951886133
$ws.Range("H22").Value = 22.17686133172641
$ws.Range("I22").Value = 34.96691491603409
$ws.Range("J22").Value = 8.065234891441255
$ws.Range("K22").Value = 15.67130696225107
$ws.Range("L22").Value = 12.97456310061315
$ws.Range("M22").Value = 18.55115099066279
$ws.Range("N22").Value = 24.31990784978969

$ws.Range("B23").Value = 18.38311842821011
$ws.Range("C23").Value = 4.900683696751009
$ws.Range("D23").Value = 13.27664647114019
$ws.Range("E23").Value = 13.04452784463926
$ws.Range("G23").Value = 60.77291676506471
$ws.Range("H23").Value = 22.18742245602531
$ws.Range("I23").Value = 34.98307457133433
$ws.Range("J23").Value = 8.064871357433795
$ws.Range("K23").Value = 15.62486124740285
$ws.Range("L23").Value = 12.9770403613171
$ws.Range("M23").Value = 18.54181941969588
$ws.Range("N23").Value = 24.3424278656554

$ws.Range("B24").Value = 18.14058546694449
$ws.Range("C24").Value = 4.584916655815572
$ws.Range("D24").Value = 13.25938133578583
$ws.Range("E24").Value = 13.07148975909595
$ws.Range("G24").Value = 60.71253463306422
$ws.Range("H24").Value = 22.23279884985824
$ws.Range("I24").Value = 35.05335411007952
$ws.Range("J24").Value = 8.063528221431175
$ws.Range("K24").Value = 15.45307637555495
$ws.Range("L24").Value = 12.98922516836636
$ws.Range("M24").Value = 18.51125508956841
$ws.Range("N24").Value = 24.43102025385899

$ws.Range("B25").Value = 17.88950980500539
$ws.Range("C25").Value = 4.215282191549387
$ws.Range("D25").Value = 13.24864754963952
$ws.Range("E25").Value = 13.10503487031436
$ws.Range("G25").Value = 60.6924379799469
$ws.Range("H25").Value = 22.29323065464114
$ws.Range("I25").Value = 35.14853490405336
$ws.Range("J25").Value = 8.062151248768636
$ws.Range("K25").Value = 15.2773367005667
$ws.Range("L25").Value = 13.00837813506054
$ws.Range("M25").Value = 18.48861018994218
$ws.Range("N25").Value = 24.53366344283263
